$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")

# Update G8 status: was "Need to do review..." -> now "Ready for printing"
$ws.Range("G8").Value = "Ready for printing"

# Drop the wrap-text formatting this status no longer needs (match style of F8)
$ws.Range("G8").WrapText = $false

# Row 8 no longer has an explicit wrapped height; let it autofit back to default
$ws.Rows.Item(8).AutoFit()

# Move the active selection to G14 (cosmetic, matches the author's last click)
$ws.Range("G14").Select()
